# Re-label the maturity column: the header "months" becomes "maturities",
# and each numeric month count in column A (1, 2, ..., 600) becomes a text
# label like "1M", "18M", "600M". "quotes" header and all rate values in
# column B are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 ("quotes") is unchanged.
$ws.Range("B1").Value = "quotes"

# Column A originally held plain numeric month counts (1..12, then 18, 24,
# 36, ..., 600). Re-write them as "<n>M" text labels first ...
$months = @(1,2,3,4,5,6,7,8,9,10,11,12,18,24,36,48,60,72,84,96,108,120,132,144,180,240,300,360,480,600)

for ($i = 0; $i -lt $months.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "$($months[$i])M"
}

# ... then relabel the header from "months" to "maturities".
$ws.Range("A1").Value = "maturities"

# The new labels are a little wider than the plain numbers were, so column A
# is widened to fit; columns B:C keep their original width.
[void]$ws.Range("A1:A31").EntireColumn.AutoFit()

# Active selection moves back to A2.
[void]$ws.Range("A2").Select()
